$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "58.542.44"
$ws.Cells.Item(2, 5).Value = "  -1.91%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "2.578.06"
$ws.Cells.Item(3, 5).Value = "  -2.59%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  +0.02%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "548.70"
$ws.Cells.Item(5, 5).Value = "  +1.95%  "

# Row 6
$ws.Cells.Item(6, 4).Value = "144.07"
$ws.Cells.Item(6, 5).Value = "  -0.95%  "

# Row 7
$ws.Cells.Item(7, 4).Value = "0.999"
$ws.Cells.Item(7, 5).Value = "  +0.06%  "

# Row 8
$ws.Cells.Item(8, 4).Value = "0.589"
$ws.Cells.Item(8, 5).Value = "  +2.59%  "

# Row 9
$ws.Cells.Item(9, 4).Value = "6.84"
$ws.Cells.Item(9, 5).Value = "  +2.97%  "

# Row 10
$ws.Cells.Item(10, 4).Value = "0.100"
$ws.Cells.Item(10, 5).Value = "  -2.78%  "

# Row 11
$ws.Cells.Item(11, 4).Value = "0.139"
$ws.Cells.Item(11, 5).Value = "  +3.73%  "

# Row 12
$ws.Cells.Item(12, 5).Value = "  -1.86%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "3.035.18"
$ws.Cells.Item(13, 5).Value = "  -2.79%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "58.496.05"
$ws.Cells.Item(14, 5).Value = "  -1.85%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "20.61"
$ws.Cells.Item(15, 5).Value = "  -2.53%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "2.579.23"
$ws.Cells.Item(16, 5).Value = "  -2.65%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "0.0000131"
$ws.Cells.Item(17, 5).Value = "  -3.21%  "

# Row 18
$ws.Cells.Item(18, 5).Value = "  +0.10%  "

# Row 19
$ws.Cells.Item(19, 4).Value = "334.28"
$ws.Cells.Item(19, 5).Value = "  -2.11%  "

# Row 20
$ws.Cells.Item(20, 5).Value = "  -4.03%  "

# Row 21
$ws.Cells.Item(21, 4).Value = "6.09"
$ws.Cells.Item(21, 5).Value = "  -3.76%  "

# Row 22
$ws.Cells.Item(22, 4).Value = "0.999"
$ws.Cells.Item(22, 5).Value = "  +0.04%  "

# Row 23
$ws.Cells.Item(23, 4).Value = "66.51"
$ws.Cells.Item(23, 5).Value = "  -0.30%  "

# Row 24
$ws.Cells.Item(24, 4).Value = "0.422"
$ws.Cells.Item(24, 5).Value = "  +1.32%  "

# Row 25
$ws.Cells.Item(25, 5).Value = "  -0.14%  "

# Row 26
$ws.Cells.Item(26, 4).Value = "0.158"
$ws.Cells.Item(26, 5).Value = "  -4.86%  "

# Row 27
$ws.Cells.Item(27, 5).Value = "  -3.87%  "

# Row 28
$ws.Cells.Item(28, 4).Value = "0.0₃0738"
$ws.Cells.Item(28, 5).Value = "  -1.98%  "

# Row 29
$ws.Cells.Item(29, 5).Value = "  +0.02%  "

# Row 30
$ws.Cells.Item(30, 4).Value = "1.66"
$ws.Cells.Item(30, 5).Value = "  +0.06%  "

# Row 31
$ws.Cells.Item(31, 4).Value = "154.99"
$ws.Cells.Item(31, 5).Value = "  +2.65%  "

# Row 32
$ws.Cells.Item(32, 4).Value = "5.88"
$ws.Cells.Item(32, 5).Value = "  +0.39%  "

# Row 33
$ws.Cells.Item(33, 4).Value = "18.80"
$ws.Cells.Item(33, 5).Value = "  -0.87%  "

# Row 34
$ws.Cells.Item(34, 4).Value = "3.89"
$ws.Cells.Item(34, 5).Value = "  -3.07%  "

# Row 35
$ws.Cells.Item(35, 4).Value = "37.16"
$ws.Cells.Item(35, 5).Value = "  -0.15%  "

# Row 36
$ws.Cells.Item(36, 4).Value = "0.851"
$ws.Cells.Item(36, 5).Value = "  +1.41%  "

# Row 37
$ws.Cells.Item(37, 4).Value = "1.10"
$ws.Cells.Item(37, 5).Value = "  -3.38%  "

# Row 38
$ws.Cells.Item(38, 2).Value = "Stacks"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(38, 4).Value = "1.43"
$ws.Cells.Item(38, 5).Value = "  -1.76%  "

# Row 39
$ws.Cells.Item(39, 2).Value = "Fetch.AI"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Cells.Item(39, 4).Value = "0.818"
$ws.Cells.Item(39, 5).Value = "  -2.52%  "

# Row 40
$ws.Cells.Item(40, 4).Value = "3.57"
$ws.Cells.Item(40, 5).Value = "  -0.67%  "

# Row 41
$ws.Cells.Item(41, 4).Value = "277.85"
$ws.Cells.Item(41, 5).Value = "  -4.39%  "

# Row 42
$ws.Cells.Item(42, 4).Value = "0.999"
$ws.Cells.Item(42, 5).Value = "  +0.03%  "

# Row 43
$ws.Cells.Item(43, 5).Value = "  -2.49%  "

# Row 44
$ws.Cells.Item(44, 4).Value = "10.62"
$ws.Cells.Item(44, 5).Value = "  -1.06%  "

# Row 45
$ws.Cells.Item(45, 4).Value = "0.0946"
$ws.Cells.Item(45, 5).Value = "  -0.21%  "

# Row 46
$ws.Cells.Item(46, 5).Value = "  -1.94%  "

# Row 47
$ws.Cells.Item(47, 5).Value = "  -0.66%  "

# Row 48
$ws.Cells.Item(48, 4).Value = "1.906.90"
$ws.Cells.Item(48, 5).Value = "  -3.62%  "

# Row 49
$ws.Cells.Item(49, 5).Value = "  -3.58%  "

# Row 50
$ws.Cells.Item(50, 4).Value = "17.69"
$ws.Cells.Item(50, 5).Value = "  -3.81%  "

# Row 51
$ws.Cells.Item(51, 4).Value = "111.75"
$ws.Cells.Item(51, 5).Value = "  +1.09%  "

